$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Unit Processes")

# --- Update the unit-library file paths to be relative to the repo root ---
# (previously "units/....xlsx", now "demoData/toy cement/units/....xlsx")
$ws.Range("D2").Value = "demoData/toy cement/units/cementUnits.xlsx"
$ws.Range("F2").Value = "demoData/toy cement/units/cementUnits.xlsx"

$ws.Range("D3").Value = "demoData/toy cement/units/cementUnits.xlsx"
$ws.Range("F3").Value = "demoData/toy cement/units/cementUnits.xlsx"

$ws.Range("D4").Value = "demoData/toy cement/units/cementUnits.xlsx"
$ws.Range("F4").Value = "demoData/toy cement/units/cementUnits.xlsx"

$ws.Range("D5").Value = "demoData/toy cement/units/auxUnits.xlsx"
$ws.Range("F5").Value = "demoData/toy cement/units/auxUnits.xlsx"

$ws.Range("D6").Value = "demoData/toy cement/units/auxUnits.xlsx"
$ws.Range("F6").Value = "demoData/toy cement/units/auxUnits.xlsx"

$ws.Range("D7").Value = "demoData/toy cement/units/auxUnits.xlsx"
$ws.Range("F7").Value = "demoData/toy cement/units/auxUnits.xlsx"

# --- Fix the swapped/mislabeled varSheet (E) and calcSheet (G) values, ---
# --- also correcting the "v CO2 capture" typo to "CO2 Capture"         ---
$ws.Range("E5").Value = "v CO2 Capture"
$ws.Range("G5").Value = "c CO2 Capture"

$ws.Range("E6").Value = "v CO2 Compression"
$ws.Range("G6").Value = "c CO2 Compression"

$ws.Range("E7").Value = "v Power Station"
$ws.Range("G7").Value = "c Power Station"

# --- Add the new "mysterious_cement_factory" unit process row ---
$ws.Range("A8").Value = "mysterious_cement_factory"
$ws.Range("A8").NumberFormat = "@"

$ws.Range("B8").Value = "cement"
$ws.Range("B8").NumberFormat = "@"

$ws.Range("C8").Value = "output"
$ws.Range("C8").NumberFormat = "@"

$ws.Range("D8").Value = "demoData/toy cement/units/cementUnits.xlsx"
$ws.Range("D8").NumberFormat = "@"

$ws.Range("E8").Value = "var mystery factory"

$ws.Range("F8").Value = "demoData/toy cement/units/cementUnits.xlsx"
$ws.Range("F8").NumberFormat = "@"

$ws.Range("G8").Value = "c mystery factory"

$ws.Range("H8").Value = "this factory is very mysterious"
$ws.Range("H8").NumberFormat = "@"

# --- Update the selection to match the author's final cursor position ---
[void]$ws.Range("G5").Select()
